# "Actualizar 02-05-2021 02-13-38" automated refresh.
# 1) The previous batch of status rows (268-281) gets its "Fecha" timestamp
#    normalised/refreshed to 44232.0716925.
# 2) A brand-new batch of 14 status rows (282-295) is appended, following
#    the same repeating 14-row cycle (Odoo/Blackbox/PowerBI/.../EZ Exporter)
#    with a fresh timestamp of 44232.09276610406.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Refresh the timestamp on the last existing batch (rows 268-281) ---
for ($r = 268; $r -le 281; $r++) {
    $ws.Cells.Item($r, 4).Value = 44232.0716925
}

# --- 2) Append the new batch (rows 282-295) ---

# Name (col A), status (col C) shared strings used by the repeating cycle.
$names = @("Odoo","Blackbox","PowerBI","Dropbox","Odoo","GEE","UtilidadesOdoo","Filtros Dashboard","MapStore","GeoServer","Tomcat","Shiny","Github","EZ Exporter")

# URL text shown in col B (what the user sees / what's stored as the cell value).
$urlText = @(
    "https://www.dataintelligence-group.com/",
    "https://serviciodashboard.azurewebsites.net/",
    "https://powerbi.microsoft.com/es-es/",
    "https://www.dropbox.com/",
    "https://dataintelligence.store/",
    "https://app-data-i.users.earthengine.app/",
    "https://odooutil.azurewebsites.net/",
    "https://filtradordashboard.azurewebsites.net/",
    "https://ide.dataintelligence-group.com/mapstore/#/",
    "https://ide.dataintelligence-group.com/geoserver/web/?0",
    "https://ide.dataintelligence-group.com/",
    "https://rpubs.com/dataintelligence/",
    "https://github.com/Sud-Austral/",
    "https://ezexporter.highviewapps.com/exports/export-profile/"
)

# Hyperlink target address (without any #fragment - that goes in SubAddress).
$urlAddr = @(
    "https://www.dataintelligence-group.com/",
    "https://serviciodashboard.azurewebsites.net/",
    "https://powerbi.microsoft.com/es-es/",
    "https://www.dropbox.com/",
    "https://dataintelligence.store/",
    "https://app-data-i.users.earthengine.app/",
    "https://odooutil.azurewebsites.net/",
    "https://filtradordashboard.azurewebsites.net/",
    "https://ide.dataintelligence-group.com/mapstore/",
    "https://ide.dataintelligence-group.com/geoserver/web/?0",
    "https://ide.dataintelligence-group.com/",
    "https://rpubs.com/dataintelligence/",
    "https://github.com/Sud-Austral/",
    "https://ezexporter.highviewapps.com/exports/export-profile/"
)

# SubAddress ("location") - only the MapStore entry (cycle index 8) carries one.
$urlSub = @($null, $null, $null, $null, $null, $null, $null, $null, "/", $null, $null, $null, $null, $null)

$newTimestamp = 44232.09276610406
$startRow = 282

for ($i = 0; $i -lt 14; $i++) {
    $row = $startRow + $i

    # Copy the formatting (styles) from the equivalent row one cycle back
    # (268 + i) so fonts/number-formats/borders match exactly.
    $srcRow = 268 + $i
    $ws.Range("A" + $srcRow + ":D" + $srcRow).Copy() | Out-Null
    $ws.Range("A" + $row + ":D" + $row).PasteSpecial(-4122) | Out-Null  # xlPasteFormats

    $ws.Cells.Item($row, 1).Value = $names[$i]
    $ws.Cells.Item($row, 2).Value = $urlText[$i]
    $ws.Cells.Item($row, 3).Value = "Disponible"
    $ws.Cells.Item($row, 4).Value = $newTimestamp

    $target = $ws.Cells.Item($row, 2)
    if ($urlSub[$i]) {
        $ws.Hyperlinks.Add($target, $urlAddr[$i], $urlSub[$i]) | Out-Null
    } else {
        $ws.Hyperlinks.Add($target, $urlAddr[$i]) | Out-Null
    }

    # Hyperlinks.Add() re-stamps a style on col B - put the normal hyperlink
    # style (copied above from the previous cycle) back on top of it.
    $ws.Range("B" + $srcRow).Copy() | Out-Null
    $ws.Range("B" + $row).PasteSpecial(-4122) | Out-Null  # xlPasteFormats
}

$excel.CutCopyMode = 0
